# Reassign the order of the statistic columns (C:F) on the active sheet.
# Previous order: C=Moda, D=Media, E=Cuasidesviación, F=Mediana
# New order:      C=Media, D=Cuasidesviación, E=Mediana, F=Moda
# (i.e. a left-rotation of the four columns' contents)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("C1").Value = "Media"
$ws.Range("D1").Value = "Cuasidesviación"
$ws.Range("E1").Value = "Mediana"
$ws.Range("F1").Value = "Moda"

# --- Libro 68 data row (row 2, merged down through row 7) ---
$ws.Range("C2").Value = 5.283236994219653
$ws.Range("D2").Value = 1.213110041412962
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 6

# --- Libro 69 data row (row 8, merged down through row 13) ---
$ws.Range("C8").Value = 5.682539682539683
$ws.Range("D8").Value = 0.8948850498428026
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 6
